$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, B, C, D, E, DneedsText
$data = @(
    ,(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "61.141.89", "  -0.96%  ", 0)
    ,(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.386.54", "  -0.01%  ", 0)
    ,(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.11%  ", 1)
    ,(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "574.44", "  -0.77%  ", 1)
    ,(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "137.27", "  +0.08%  ", 1)
    ,(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.02%  ", 1)
    ,(8, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.383.58", "  -0.08%  ", 0)
    ,(9, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.470", "  -1.11%  ", 1)
    ,(10, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "7.64", "  +2.17%  ", 1)
    ,(11, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.122", "  -2.86%  ", 1)
    ,(12, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.381", "  -2.43%  ", 1)
    ,(13, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.961.41", "  -0.26%  ", 0)
    ,(14, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.123", "  +0.65%  ", 1)
    ,(15, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000172", "  -2.54%  ", 1)
    ,(16, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "25.73", "  +1.49%  ", 1)
    ,(17, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.380.44", "  -0.46%  ", 0)
    ,(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "61.261.45", "  -0.93%  ", 0)
    ,(19, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "13.84", "  -1.94%  ", 1)
    ,(20, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.76", "  -0.86%  ", 1)
    ,(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "9.35", "  -1.41%  ", 1)
    ,(22, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "377.12", "  -0.85%  ", 1)
    ,(23, "WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "3.519.00", "  -0.35%  ", 0)
    ,(24, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.551", "  -2.07%  ", 1)
    ,(25, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.07%  ", 1)
    ,(26, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0000126", "  -0.32%  ", 1)
    ,(27, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "71.00", "  -0.32%  ", 1)
    ,(28, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.179", "  +11.80%  ", 1)
    ,(29, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.65", "  -2.84%  ", 1)
    ,(30, "Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "1.00", "  +0.04%  ", 1)
    ,(31, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "7.44", "  -2.19%  ", 1)
    ,(32, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "8.08", "  -1.57%  ", 1)
    ,(33, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.15", "  -1.67%  ", 1)
    ,(34, "USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  -0.05%  ", 1)
    ,(35, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "23.44", "  -0.05%  ", 1)
    ,(36, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "5.15", "  -4.23%  ", 1)
    ,(37, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.55", "  -1.53%  ", 1)
    ,(38, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.84", "  -0.55%  ", 1)
    ,(39, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "164.67", "  -0.25%  ", 1)
    ,(40, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0761", "  -3.30%  ", 1)
    ,(41, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.999", "  -0.10%  ", 1)
    ,(42, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "25.47", "  +2.31%  ", 1)
    ,(43, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.776", "  -1.03%  ", 1)
    ,(44, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.70", "  -1.30%  ", 1)
    ,(45, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.35", "  -1.59%  ", 1)
    ,(46, "ONDO", "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo", "1.19", "  -3.56%  ", 1)
    ,(47, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.554.90", "  +8.92%  ", 0)
    ,(48, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "6.79", "  -1.15%  ", 1)
    ,(49, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "22.95", "  +0.05%  ", 1)
    ,(50, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.45", "  +4.38%  ", 1)
    ,(51, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0259", "  -1.27%  ", 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    if ($row[5] -eq 1) {
        $ws.Range("D$r").NumberFormat = "@"
    }
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}
